$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata") updates ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting everything below it up by one row (dimension A1:B21 -> A1:B20)
$ws1.Rows.Item(11).Delete()

# --- Sheet 2 ("Elements") updates ---
$ws2 = $wb.Worksheets.Item(2)

# Top-level Extension row: Short/Definition columns (K2/L2) updated from the generic
# "Extension" / "An Extension" text to the profile-specific title/description
$ws2.Range("K2").Value = "Region Code"
$ws2.Range("L2").Value = "Customer-specific code for the geographic region of the address"
